$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells (registry value reassignments) ---
# Order matters for shared-string table layout: match the authoring sequence.
$ws.Range("B21").Value = "T_MSGHASH"
$ws.Range("B11").Value = "unassigned"
$ws.Range("F11").Value = "unassigned"
$ws.Range("F13").Value = "unassigned"
$ws.Range("B22").Value = "unassigned"
$ws.Range("B23").Value = "unassigned"
$ws.Range("F16").Value = "T_RSA-SHA256"

# --- Add new "HashFunctions" rows 27-29 ---
$ws.Range("A27").Value = "HashFunctions"
$ws.Range("B27").Value = "T_SHA-256"
$ws.Range("C27").Value = "%x0001"
$ws.Range("D27").Value = "messages"

$ws.Range("A28").Value = "HashFunctions"
$ws.Range("B28").Value = "T_SHA-512"
$ws.Range("C28").Value = "%x0002"
$ws.Range("D28").Value = "messages"

$ws.Range("A29").Value = "HashFunctions"
$ws.Range("C29").Value = "%x1000 - %x1FFF"
$ws.Range("B29").Value = "reserved"
$ws.Range("D29").Value = "messages"

# Match formatting of column A "level" cells in the new rows
$ws.Range("A19").Copy()
$ws.Range("A27:A29").PasteSpecial(-4122)

# Column F width (new content added to column F)
$ws.Columns.Item(6).ColumnWidth = 26.83

# View state
$ws.Range("F26").Select()
